$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -3
    "F5"  = 2
    "F6"  = 1
    "F9"  = 3
    "F11" = -7
    "F13" = 0
    "F14" = -4
    "F15" = -1
    "F16" = 3
    "F17" = -3
    "F18" = 6
    "F19" = 2
    "F20" = -4
    "F21" = -2
    "F22" = -3
    "F23" = 2
    "F24" = -5
    "F25" = 3
    "F26" = 4
    "F27" = 3
    "F28" = 0
    "F29" = -5
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
